$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells I1 and J1, copying the formatting from H1 (the last
# existing header cell) so they share the same bold/border/centered style.
$ws.Cells.Item(1, 9).Value = "I0"
$ws.Cells.Item(1, 10).Value = "IF"
$ws.Cells.Item(1, 8).Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Per-row I0 / IF values (column I = I0, column J = IF) keyed by sheet row number.
$data = @{
    2 = @(5, 6)
    3 = @(5, 6)
    4 = @(6, 6)
    5 = @(8, 8)
    6 = @(7, 7)
    7 = @(7, 7)
    8 = @(8, 8)
    9 = @(7, 7)
    10 = @(6, 6)
    11 = @(8, 8)
    12 = @(5, 5)
    13 = @(7, 7)
    14 = @(7, 7)
    15 = @(7, 7)
    16 = @(11, 11)
    17 = @(6, 7)
    18 = @(5, 5)
    19 = @(7, 7)
    20 = @(7, 7)
    21 = @(8, 8)
    22 = @(9, 9)
    23 = @(8, 8)
    24 = @(9, 9)
    25 = @(7, 7)
    26 = @(8, 8)
    27 = @(7, 7)
    28 = @(7, 7)
    29 = @(7, 8)
    30 = @(7, 7)
    31 = @(7, 7)
    32 = @(7, 7)
    33 = @(8, 8)
    34 = @(8, 8)
    35 = @(7, 7)
    36 = @(8, 8)
    37 = @(6, 7)
    38 = @(8, 8)
    39 = @(7, 8)
    40 = @(8, 8)
    41 = @(7, 8)
    42 = @(7, 7)
    43 = @(8, 8)
    44 = @(7, 7)
    45 = @(7, 7)
    46 = @(8, 9)
    47 = @(8, 8)
    48 = @(7, 7)
    49 = @(7, 7)
    50 = @(8, 8)
    51 = @(6, 7)
    52 = @(7, 7)
    53 = @(7, 7)
    54 = @(7, 7)
    55 = @(9, 9)
    56 = @(6, 7)
    57 = @(5, 5)
    58 = @(8, 8)
    59 = @(7, 7)
    60 = @(6, 6)
    61 = @(7, 8)
    62 = @(7, 7)
    63 = @(6, 6)
    64 = @(6, 6)
    65 = @(5, 5)
    66 = @(7, 7)
    67 = @(9, 9)
    68 = @(7, 8)
    69 = @(8, 8)
    70 = @(4, 4)
    71 = @(5, 5)
    72 = @(3, 3)
    73 = @(6, 6)
    74 = @(7, 7)
}

foreach ($key in $data.Keys) {
    $row = [int]$key
    $vals = $data[$key]
    $ws.Cells.Item($row, 9).Value = $vals[0]
    $ws.Cells.Item($row, 10).Value = $vals[1]
}
